# formulaires.xlsx - jxls template placeholder syntax migration
#
# The commit switches the jxls expression-language delimiters used in the
# template cells from the legacy "$(...)" form to the modern "${...}" form,
# and also leaves the selected/active cell on C5 instead of H7.
#
# Cell layout (row 2 = formulaire header, rows 3/4 = section rows):
#   B2 -> ${formulaire.name}
#   H2 -> ${formulaire.numero}
#   B3 -> ${section.ordre}
#   C3 -> ${section.name}
#   B4 -> ${section.ordre}
#   C4 -> ${section.name}
#
# Cell formatting (borders/fills/alignment) is left untouched - only the
# text content of the placeholder cells and the sheet's active selection
# change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "`${formulaire.name}"
$ws.Range("H2").Value = "`${formulaire.numero}"
$ws.Range("B3").Value = "`${section.ordre}"
$ws.Range("C3").Value = "`${section.name}"
$ws.Range("B4").Value = "`${section.ordre}"
$ws.Range("C4").Value = "`${section.name}"

# Move the active selection from H7 to C5, matching the authored sheetView.
$ws.Range("C5").Select()
